$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values
$ws.Range("D2").Value = "64.010.62"
$ws.Range("D3").Value = "3.514.06"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "3.514.39"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.486"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Value = "4.106.82"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Value = "3.507.49"
$ws.Range("D18").Value = "64.074.67"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "383.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Value = "3.653.08"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "3.522.97"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "160.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.811"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("D50").Value = "2.472.74"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.79"
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) (column E) values
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("E12").Value = "  -2.44%  "
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("E19").Value = "  -3.02%  "
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("E23").Value = "  -1.16%  "
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("E28").Value = "  +2.79%  "
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("E30").Value = "  -2.38%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("E33").Value = "  -1.85%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  -2.31%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("E41").Value = "  -4.40%  "
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("E43").Value = "  +3.75%  "
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -2.99%  "
$ws.Range("E47").Value = "  -3.66%  "
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("E49").Value = "  -3.14%  "
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("E51").Value = "  -1.62%  "
